# Insert a new weekly data row at row 105 (pushing existing rows 105-219
# down to 106-220), matching the commit "Fruta / hortaliza, semanal".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(105).Insert()

$ws.Range("A105").Value = 3
$ws.Range("B105").Value = "Femacal de La Calera"
$ws.Range("C105").Value = "Coquimbo"
$ws.Range("D105").Value = 44494
$ws.Range("E105").Value = 5
$ws.Range("F105").Value = 100112043
$ws.Range("G105").Value = "Pepino ensalada"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 105
$ws.Range("K105").Value = 6500
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = 6738
$ws.Range("N105").Value = "$/caja 70 unidades"
$ws.Range("O105").Value = "Región de Arica y Parinacota"
$ws.Range("P105").Value = 96
$ws.Range("Q105").Value = 70
$ws.Range("R105").Value = "Hortaliza"
